$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2159.0386
$ws.Range("I17").Value = 800
$ws.Range("J17").Value = 2336.3044
$ws.Range("K17").Value = 2400
$ws.Range("L17").Value = 7008.9132
$ws.Range("M17").Value = -2232
$ws.Range("N17").Value = -7344.9132

$ws.Range("H86").Value = 3142.4285
$ws.Range("J86").Value = 3166.3333
$ws.Range("L86").Value = 3166.3333
$ws.Range("N86").Value = -5412.3333

$ws.Range("H89").Value = 3142.4285
$ws.Range("J89").Value = 3166.3333
$ws.Range("L89").Value = 15831.6665
$ws.Range("N89").Value = -27063.6665

$ws.Range("H99").Value = 138
$ws.Range("I99").Value = 138
$ws.Range("K99").Value = 414
$ws.Range("M99").Value = 1084

$ws.Range("H112").Value = 2110.5454
$ws.Range("J112").Value = 1673.7142
$ws.Range("L112").Value = 5021.142599999999
$ws.Range("N112").Value = -7237.142599999999

$ws.Range("H135").Value = 1157.091
$ws.Range("I135").Value = 972.8
$ws.Range("J135").Value = 3000
$ws.Range("K135").Value = 8755.199999999999
$ws.Range("L135").Value = 27000
$ws.Range("M135").Value = -6220.199999999999
$ws.Range("N135").Value = -32070

$ws.Range("H137").Value = 2112.524
$ws.Range("I137").Value = 1232.1111
$ws.Range("K137").Value = 3696.3333
$ws.Range("M137").Value = -1146.3333

$ws.Range("H138").Value = 2514.3333
$ws.Range("I138").Value = 681.3333
$ws.Range("J138").Value = 4958.3335
$ws.Range("K138").Value = 2043.9999
$ws.Range("L138").Value = 14875.0005
$ws.Range("M138").Value = 3096.0001
$ws.Range("N138").Value = -25155.0005

$ws.Range("H141").Value = 1528.0625
$ws.Range("I141").Value = 1419.3077
$ws.Range("J141").Value = 1999.3334
$ws.Range("K141").Value = 4257.9231
$ws.Range("L141").Value = 5998.0002
$ws.Range("M141").Value = 922.0769
$ws.Range("N141").Value = -16358.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 11572.05
$ws.Range("J44").Value = 11865.315
$ws.Range("L44").Value = 11865.315
$ws.Range("N44").Value = -12841.315

$ws.Range("H61").Value = 2512.5557
$ws.Range("I61").Value = 1659.8572
$ws.Range("K61").Value = 1659.8572
$ws.Range("M61").Value = -1447.8572

$ws.Range("H74").Value = 3350.3333
$ws.Range("I74").Value = 2624
$ws.Range("K74").Value = 2624
$ws.Range("M74").Value = -1750

$ws.Range("H77").Value = 3350.3333
$ws.Range("I77").Value = 2624
$ws.Range("K77").Value = 13120
$ws.Range("M77").Value = -8752

$ws.Range("H102").Value = 17860016
$ws.Range("I102").Value = 20835202
$ws.Range("J102").Value = 8899
$ws.Range("K102").Value = 20835202
$ws.Range("L102").Value = 8899
$ws.Range("M102").Value = -20833580
$ws.Range("N102").Value = -12143

$ws.Range("H132").Value = 823.0714
$ws.Range("I132").Value = 833.92
$ws.Range("J132").Value = 732.6667
$ws.Range("K132").Value = 2501.76
$ws.Range("L132").Value = 2198.0001
$ws.Range("M132").Value = 28.24000000000024
$ws.Range("N132").Value = -7258.0001

$ws.Range("H136").Value = 2512.5557
$ws.Range("I136").Value = 1659.8572
$ws.Range("K136").Value = 4979.571599999999
$ws.Range("M136").Value = -2429.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 79
$ws.Range("I5").Value = 116
$ws.Range("K5").Value = 116
$ws.Range("M5").Value = -3

$ws.Range("H134").Value = 5116.0713
$ws.Range("I134").Value = 974.6
$ws.Range("J134").Value = 15469.75
$ws.Range("K134").Value = 2923.8
$ws.Range("L134").Value = 46409.25
$ws.Range("M134").Value = -388.8000000000002
$ws.Range("N134").Value = -51479.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4888.6924
$ws.Range("I31").Value = 1763.5
$ws.Range("K31").Value = 1763.5
$ws.Range("M31").Value = -1468.5

$ws.Range("H34").Value = 4888.6924
$ws.Range("I34").Value = 1763.5
$ws.Range("K34").Value = 1763.5
$ws.Range("M34").Value = -1561.5

$ws.Range("H58").Value = 3759
$ws.Range("I58").Value = 2450
$ws.Range("K58").Value = 2450
$ws.Range("M58").Value = -2247

$ws.Range("H136").Value = 3759
$ws.Range("I136").Value = 2450
$ws.Range("K136").Value = 7350
$ws.Range("M136").Value = -4800

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1966.5385
$ws.Range("J34").Value = 3535.1428
$ws.Range("L34").Value = 10605.4284
$ws.Range("N34").Value = -10773.4284

$ws.Range("H132").Value = 3265.6428
$ws.Range("I132").Value = 2496.6667
$ws.Range("K132").Value = 22470.0003
$ws.Range("M132").Value = -19940.0003

$ws.Range("H137").Value = 4716.6665
$ws.Range("I137").Value = 2600
$ws.Range("K137").Value = 7800
$ws.Range("M137").Value = -2700

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 3000
$ws.Range("I46").Value = 3000
$ws.Range("K46").Value = 3000
$ws.Range("M46").Value = -2844

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()

$ws.Range("H102").Value = 2403.5264
$ws.Range("I102").Value = 1476.2858
$ws.Range("K102").Value = 1476.2858
$ws.Range("M102").Value = 145.7141999999999

$ws.Range("H113").Value = 3251.5454
$ws.Range("I113").Value = 2171
$ws.Range("J113").Value = 6133
$ws.Range("K113").Value = 2171
$ws.Range("L113").Value = 6133
$ws.Range("M113").Value = -1
$ws.Range("N113").Value = -10473

$ws.Range("H122").Value = 2412
$ws.Range("I122").Value = 1947
$ws.Range("K122").Value = 5841
$ws.Range("M122").Value = -3391

$ws.Range("H126").Value = 2999.6667
$ws.Range("I126").Value = 2499.5
$ws.Range("K126").Value = 7498.5
$ws.Range("M126").Value = -5028.5

$ws.Range("H132").Value = 1760.7368
$ws.Range("I132").Value = 1760.7368
$ws.Range("K132").Value = 5282.2104
$ws.Range("M132").Value = -2752.2104

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 7500
$ws.Range("I53").Value = 5000
$ws.Range("K53").Value = 5000
$ws.Range("M53").Value = -4482

$ws.Range("H122").Value = 2382
$ws.Range("I122").Value = 2527.9167
$ws.Range("K122").Value = 7583.750100000001
$ws.Range("M122").Value = -5133.750100000001

$ws.Range("H132").Value = 3383.7144
$ws.Range("I132").Value = 3538.2
$ws.Range("K132").Value = 10614.6
$ws.Range("M132").Value = -8084.599999999999

$ws.Range("H136").Value = 3666.6667
$ws.Range("I136").Value = 3666.6667
$ws.Range("K136").Value = 11000.0001
$ws.Range("M136").Value = -8450.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 12610
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H52").Value = 24499.666
$ws.Range("I52").Value = 34999.5
$ws.Range("J52").Value = 3500
$ws.Range("K52").Value = 34999.5
$ws.Range("L52").Value = 3500
$ws.Range("M52").Value = -34773.5
$ws.Range("N52").Value = -3952

$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -59178

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H136").Value = 3378.2666
$ws.Range("I136").Value = 2315.111
$ws.Range("K136").Value = 6945.333
$ws.Range("M136").Value = -4395.333

Write-Output "edits applied"